# BOM update for 1st SETI manufacturing batch
$wb  = $excel.ActiveWorkbook
$ws  = $wb.Worksheets.Item("BATT_CPU")
$wsh = $wb.Worksheets.Item("_HISTORY")

# --- Capacitor values: add voltage & dielectric info ---
$ws.Range("E8").Value  = "2.2uF, 16V, X7R"
$ws.Range("E9").Value  = "1.0uF, 16V, X7R"
$ws.Range("E10").Value = "22pF, 16V, C0G"
$ws.Range("E11").Value = "4.7uF, 16V, X7R"
$ws.Range("E12").Value = "100nF, 25V, X7R o X5R"

# --- L7: changed from inductor to 0R resistor ---
$ws.Range("E13").Value = "0R"
$ws.Range("F13").Value = "RESISTOR"

# --- pos 8 (Q1,Q2,Q5): package corrected ---
$ws.Range("G14").Value = "SOT323"

# --- pos 9 (Q3,Q4): new datasheet link ---
$ws.Range("H15").Value = "https://www.digikey.es/es/products/detail/nexperia-usa-inc/2N7002PW-115/2296328?s=N4IgTCBcDa4HYHYAMSwAcDuIC6BfIA"
$ws.Hyperlinks.Add($ws.Range("H15"), "https://www.digikey.es/es/products/detail/nexperia-usa-inc/2N7002PW-115/2296328?s=N4IgTCBcDa4HYHYAMSwAcDuIC6BfIA") | Out-Null
$ws.Range("H15").Style = $ws.Range("H14").Style

# --- pos 15: add R5 to the group, bump qty 4 -> 5 ---
$ws.Range("C21").Value = 5
$ws.Range("D21").Value = "R4, R5, R14, R16, R18"

# --- pos 17 (R5): deleted, merged into pos 15 ---
$ws.Range("D23").Value = "deleted"
$ws.Range("E23").Value = "deleted"
$ws.Range("F23").ClearContents()
$ws.Range("G23").ClearContents()

# --- pos 23 (U6): package corrected + new link ---
$ws.Range("G29").Value = "SOT-753"
$ws.Range("H29").Value = "https://www.digikey.es/es/products/detail/analog-devices-inc-maxim-integrated/MAX9062EUK-T/1937837"
$ws.Hyperlinks.Add($ws.Range("H29"), "https://www.digikey.es/es/products/detail/analog-devices-inc-maxim-integrated/MAX9062EUK-T/1937837") | Out-Null
$ws.Range("H29").Style = $ws.Range("H14").Style

# --- pos 24 (U8): package corrected + new link ---
$ws.Range("G30").Value = "SOT23-6"
$ws.Range("H30").Value = "https://www.digikey.es/es/products/detail/texas-instruments/SN74LVC1G19DBVR/654738?s=N4IgTCBcDaIMoDkDsAWAMgNQMIEYDiOAnACIBCGASiALoC%2BQA"
$ws.Hyperlinks.Add($ws.Range("H30"), "https://www.digikey.es/es/products/detail/texas-instruments/SN74LVC1G19DBVR/654738?s=N4IgTCBcDaIMoDkDsAWAMgNQMIEYDiOAnACIBCGASiALoC%2BQA") | Out-Null
$ws.Range("H30").Style = $ws.Range("H14").Style

# --- pos 1 (BZ1): new datasheet link ---
$ws.Range("H7").Value = "https://www.digikey.es/es/products/detail/murata-electronics/PKMCS0909E4000-R1/4878400?s=N4IgTCBcDaIAoGkCyBhAygBgJzYKIBYMiBaAJQEYQBdAXyA"
$ws.Hyperlinks.Add($ws.Range("H7"), "https://www.digikey.es/es/products/detail/murata-electronics/PKMCS0909E4000-R1/4878400?s=N4IgTCBcDaIAoGkCyBhAygBgJzYKIBYMiBaAJQEYQBdAXyA") | Out-Null
$ws.Range("H7").Style = $ws.Range("H14").Style

# --- _HISTORY: new row 4 ---
$wsh.Range("A7").Value = 4
$wsh.Range("B7").Value = "5-ene-2023"
$wsh.Range("C7").Value = "DGB"
$wsh.Range("D7").Value = "Se cambia L7 a resistencia 0R. Se añade voltaje y dielectrico de los condensadores. Se corrige encapsulado de pos14. Borrada pos17 se añade 1 a pos15"

# --- BOM Version now computed from history ---
$ws.Range("C2").Formula = "=MAX(_HISTORY!A4:A43)"

# --- view-state touch ups (zoom + remembered selection) ---
$ws.Activate()
$excel.ActiveWindow.Zoom = 120

$wsh.Activate()
$wsh.Range("C13").Select()

$ws.Activate()
$ws.Range("C3").Select()
